$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: new procedure row (obtenerCantidadDeBonosDisponiblesPorAfiliado)
$ws.Range("A38").Value = "obtenerCantidadDeBonosDisponiblesPorAfiliado(idAfiliado)"
$ws.Range("B38").Value = "cantidadBonos"

# Row 27: registrarAtencionMedica(...) -- new info text, marked Done ("Si")
$ws.Range("C27").Value = "Actualizar consulta médica con resultados de la consulta."
$ws.Range("D27").Value = "Si"

# Row 26: registrarLlegada(...) -- info text extended to mention bono, and marked Done ("Si")
$ws.Range("C26").Value = "Crear una consulta médica para ese turno y asociar un bono a esa consulta."
$ws.Range("D26").Value = "Si"

# Row 28: obtenerTurnosDeAfiliado(idAfiliado) -- Return column extended with extra fields, marked Done ("Si")
$ws.Range("B28").Value = "numeroTurno, fechaHorarioAtencion, nombreProfesional, apellidoProfesional, especialidad"
$ws.Range("D28").Value = "Si"

# Column width / layout tweaks
$ws.Columns.Item(1).ColumnWidth = 56.375
$ws.Columns.Item(2).ColumnWidth = 42.5
$ws.Columns.Item(3).ColumnWidth = 39.25
$ws.Columns.Item(4).ColumnWidth = 32.75
$ws.Columns.Item(8).ColumnWidth = 41.875

# Update view/selection to match author's final cursor position
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("D35").Select()
